$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.2859509529444935, 0.01388827069144083, 0.006621528401716747, -0.2238840684928703, 0.5325075328486121)
    3  = @(0.5107610139011399, 0.03137483124467404, 0.1766367752656752, -0.2646734840316442, 0.5230118728062177)
    4  = @(0.6184905866032928, 0.0518624791402557, 0.2613589834390402, -0.3990688081024865, 0.4883319378332514)
    5  = @(0.04743822908319877, -0.006300419179037398, -0.002192891532154102, -0.2062872690305563, 0.5461301469193761)
    6  = @(0.09285964204334098, 0.007732799577114901, 0.161427200727865, -0.2345201207183042, 0.5362671535811375)
    7  = @(0.4255621319500999, 0.03268483639646356, 0.2553335297719969, -0.3487914831783068, 0.48340775868858)
    8  = @(0.02601404063028697, -0.003093420118658873, 0.04203266763915399, -0.1681493645975394, 0.5582776633267116)
    9  = @(0.08358269619876381, 0.0100450360313335, 0.2362701805105318, -0.1894044383603727, 0.5521213756404834)
    10 = @(0.2650197135655772, 0.0291622248711331, 0.2795450273609112, -0.3037864330808801, 0.5110144147841393)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 3).Value = $values[0]
    $ws.Cells.Item($row, 4).Value = $values[1]
    $ws.Cells.Item($row, 5).Value = $values[2]
    $ws.Cells.Item($row, 6).Value = $values[3]
    $ws.Cells.Item($row, 7).Value = $values[4]
}
